# Update the CBC test row to reference "CH" instead of "CBC".
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A2").Value = "VA.MHV.PHR.chTest"
$ws.Range("B2").Value = "VA MHV PHR CH labTest"
